$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = "release/8.0.11"
$ws.Range("B14").Value = "X"
$ws.Range("C14").Value = "X"
$ws.Range("D14").Value = "X"
$ws.Range("E14").Value = "X"
